# ==========================================================================
# Applies the January 2024 data-refresh edit described by the commit
# "Add files via upload" to data202401.xlsx:
#   - Sheet "部门情况202401": append a new department row (普惠业务二部)
#   - Sheet "对公业务台账202401": append a new "户均" (per-account average)
#     column U
#   - Sheet "对公产品台账202401": append "户均" / "平均利率" columns J, K
#   - Sheet "个人经营贷202401": drop six discontinued product rows, append
#     "个人户均" / "逾期金额" / "不良金额" columns K, L, M, and refresh the
#     Total row
# ==========================================================================

$wb = $excel.ActiveWorkbook

# Helper: copies the number format / font / border (but not the value) of
# $srcAddr onto $dstAddr on worksheet $sheet, mirroring how a user would
# fill a new header cell by copying an existing one.
function Copy-HeaderFormat($sheet, $srcAddr, $dstAddr) {
    $sheet.Range($srcAddr).Copy()
    $sheet.Range($dstAddr).PasteSpecial(-4122) # xlPasteFormats
}

# --------------------------------------------------------------------
# Sheet 1: 部门情况202401  -- add row 10 "普惠业务二部"
# --------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("部门情况202401")

$row10 = [ordered]@{
    "A" = "普惠业务二部"
    "B" = "'0.00"
    "C" = "'0.00"
    "D" = "'0.00"
    "E" = "'0.00"
    "F" = "'280.00"
    "G" = "'2.00"
    "H" = "'5.82"
    "I" = "'280.00"
    "J" = "'0.00"
    "K" = "'0.00"
    "L" = "'0.00"
    "M" = "'0.00"
    "N" = "'0.00"
    "O" = "'0.00"
}
foreach ($col in $row10.Keys) {
    $ws1.Range($col + "10").Value = $row10[$col]
}

# --------------------------------------------------------------------
# Sheet 3: 对公业务台账202401  -- add column U "户均"
# --------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("对公业务台账202401")

Copy-HeaderFormat $ws3 "T1" "U1"
$ws3.Range("U1").Value = "户均"

$colU = [ordered]@{
    2  = "'8289.72"
    3  = "'10936.67"
    4  = "'380.09"
    5  = "'263.58"
    6  = "'97.82"
    7  = "'591.19"
    8  = "'18.81"
    9  = "'607.56"
    10 = "'0.33"
    11 = "'32968.00"
    12 = "'54153.77"
}
foreach ($r in $colU.Keys) {
    $ws3.Range("U" + $r).Value = $colU[$r]
}

# --------------------------------------------------------------------
# Sheet 4: 对公产品台账202401  -- add columns J "户均", K "平均利率"
# --------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("对公产品台账202401")

Copy-HeaderFormat $ws4 "I1" "J1"
Copy-HeaderFormat $ws4 "I1" "K1"
$ws4.Range("J1").Value = "户均"
$ws4.Range("K1").Value = "平均利率"

$colJK = @{
    2  = @(0, 0)
    3  = @(10936.67, 5.5)
    4  = @(380.09, 6.6)
    5  = @(263.58, 6.6)
    6  = @(97.81999999999999, 5.4)
    7  = @(591.1900000000001, 6.51)
    8  = @(18.81, 6.3)
    9  = @(791.67, 3.95)
    10 = @(160.76, 5)
    11 = @(800, 5.15)
    12 = @(907.27, 5.15)
    13 = @(489.94, 5)
    14 = @(500, 6.5)
    15 = @(0.65, 24)
    16 = @(0.01, 24)
    17 = @(32968, 7.5)
    18 = @(48906.46, 123.16)
}
foreach ($r in $colJK.Keys) {
    $ws4.Range("J" + $r).Value = $colJK[$r][0]
    $ws4.Range("K" + $r).Value = $colJK[$r][1]
}

# --------------------------------------------------------------------
# Sheet 5: 个人经营贷202401
#   -- drop 6 discontinued product rows
#   -- add columns K "个人户均", L "逾期金额", M "不良金额"
#   -- refresh the Total row
# --------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("个人经营贷202401")

# Delete bottom-up so earlier row numbers stay valid while iterating.
$ws5.Rows.Item(18).Delete() # 车商贷
$ws5.Rows.Item(10).Delete() # 瑞e惠-携程生意人贷
$ws5.Rows.Item(9).Delete()  # 瑞e惠-拍拍贷
$ws5.Rows.Item(4).Delete()  # 创业经营贷
$ws5.Rows.Item(3).Delete()  # 出行贷-狮桥车主经营
$ws5.Rows.Item(2).Delete()  # 出行贷-新生经营

Copy-HeaderFormat $ws5 "J1" "K1"
Copy-HeaderFormat $ws5 "J1" "L1"
Copy-HeaderFormat $ws5 "J1" "M1"
$ws5.Range("K1").Value = "个人户均"
$ws5.Range("L1").Value = "逾期金额"
$ws5.Range("M1").Value = "不良金额"

$colKLM = @{
    2  = @(173.93, 0, 0)
    3  = @(245.17, 1525, 560)
    4  = @(17.32, 0, 0)
    5  = @(37.42, 1497.23, 0)
    6  = @(53.7, 0, 0)
    7  = @(66.58, 0, 0)
    8  = @(15.24, 4992.67, 1822.64)
    9  = @(55.61, 0, 0)
    10 = @(11.95, 16.78, 0)
    11 = @(67.48, 0, 0)
    12 = @(17.21, 17.21, 17.21)
    13 = @(761.6100000000002, 8048.889999999999, 2399.85)
}
foreach ($r in $colKLM.Keys) {
    $ws5.Range("K" + $r).Value = $colKLM[$r][0]
    $ws5.Range("L" + $r).Value = $colKLM[$r][1]
    $ws5.Range("M" + $r).Value = $colKLM[$r][2]
}

# Refresh the Total row (B..J) now that six rows were removed.
$total13 = [ordered]@{
    "B" = 9570
    "C" = 386824.48
    "D" = 11417
    "E" = 436413.17
    "F" = 94.11
    "G" = 1581
    "H" = 67431.02
    "I" = 1926
    "J" = 1392
}
foreach ($col in $total13.Keys) {
    $ws5.Range($col + "13").Value = $total13[$col]
}
